# "Fix detection Value formula"
# Adds a new "Update Customer Table" source/output pair of workflow steps
# (with their Parameters rows) to the Steps and Parameters sheets, and
# re-orders the existing rows to match the new canonical ordering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Steps": rewrite rows 2-9 in the exact order the target file
# has them (this also controls the order new shared strings are
# allocated in, so it must mirror the target row/column order: row by
# row, A,B,C,D then I,J).
# ---------------------------------------------------------------------
$steps = $wb.Worksheets.Item("Steps")

$steps.Range("A2").Value = "Update_Nightly"
$steps.Range("B2").Value = "Update Customer Table Source"
$steps.Range("C2").Value = "Task"
$steps.Range("D2").Value = "[Commons] Evaluate Formula"

$steps.Range("A3").Value = "Update_Nightly"
$steps.Range("B3").Value = "Update Analysis table"
$steps.Range("C3").Value = "Task"
$steps.Range("D3").Value = "[Core] UpdateFormulaForAllInstances"
$steps.Range("I3").Value = "Success"
$steps.Range("J3").Value = "Update Customer Table Source"

$steps.Range("A4").Value = "Update_Nightly"
$steps.Range("B4").Value = "Back Testing"
$steps.Range("C4").Value = "Task"
$steps.Range("D4").Value = "[Core] UpdateFormulaForAllInstances"
$steps.Range("I4").Value = "Success"
$steps.Range("J4").Value = "Update Analysis table"

$steps.Range("A5").Value = "Update_Nightly"
$steps.Range("B5").Value = "Update variable"
$steps.Range("C5").Value = "Task"
$steps.Range("D5").Value = "[Core] UpdateFormulaForAllInstances"
$steps.Range("I5").Value = "Success"
$steps.Range("J5").Value = "Update Analysis table"

$steps.Range("A6").Value = "Update_Nightly"
$steps.Range("B6").Value = "Update Workspace Plane Status"
$steps.Range("C6").Value = "Task"
$steps.Range("D6").Value = "[Core] UpdateFormulaForAllInstances"

$steps.Range("A7").Value = "Update_Nightly"
$steps.Range("B7").Value = "Update Workspace Plane Table"
$steps.Range("C7").Value = "Task"
$steps.Range("D7").Value = "[Commons] Evaluate Formula"
$steps.Range("I7").Value = "Success"
$steps.Range("J7").Value = "Update Analysis table"

$steps.Range("A8").Value = "Update_Nightly"
$steps.Range("B8").Value = "Update Output table"
$steps.Range("C8").Value = "Task"
$steps.Range("D8").Value = "[Core] UpdateFormulaForAllInstances"
$steps.Range("I8").Value = "Success"
$steps.Range("J8").Value = "Update Analysis table"

$steps.Range("A9").Value = "Update_Nightly"
$steps.Range("B9").Value = "Update Customer Table Output"
$steps.Range("C9").Value = "Task"
$steps.Range("D9").Value = "[Commons] Evaluate Formula"
$steps.Range("I9").Value = "Success"
$steps.Range("J9").Value = "Update Output table"

# Drop the per-cell style overrides on the data rows (A2:D9) so the
# cells fall back to the default "Normal" style, matching the target
# file (no more s="24"/"23" on these cells).
$steps.Range("A2:D9").Style = "Normal"

$steps.Activate()
$steps.Range("A2:N9").Select()

# ---------------------------------------------------------------------
# Sheet "Parameters": rewrite rows 2-14 in the exact target order.
# ---------------------------------------------------------------------
$params = $wb.Worksheets.Item("Parameters")

$params.Range("A2").Value = "Update_Nightly"
$params.Range("B2").Value = "Update Customer Table Source"
$params.Range("C2").Value = "formula"
$params.Range("D2").Value = "return LIB_EWS.UpdateCustomerUnitSource();"

$params.Range("A3").Value = "Update_Nightly"
$params.Range("B3").Value = "Update Analysis table"
$params.Range("C3").Value = "Entity Type Name"
$params.Range("D3").Value = "Analysis_Unit"

$params.Range("A4").Value = "Update_Nightly"
$params.Range("B4").Value = "Update Analysis table"
$params.Range("C4").Value = "Attribute Name"
$params.Range("D4").Value = "LastUpdateAnalysisUnit"

$params.Range("A5").Value = "Update_Nightly"
$params.Range("B5").Value = "Back Testing"
$params.Range("C5").Value = "Entity Type Name"
$params.Range("D5").Value = "Analysis_Unit"

$params.Range("A6").Value = "Update_Nightly"
$params.Range("B6").Value = "Back Testing"
$params.Range("C6").Value = "Attribute Name"
$params.Range("D6").Value = "LastUpdateBackTesting"

$params.Range("A7").Value = "Update_Nightly"
$params.Range("B7").Value = "Update variable"
$params.Range("C7").Value = "Entity Type Name"
$params.Range("D7").Value = "Analysis_Unit"

$params.Range("A8").Value = "Update_Nightly"
$params.Range("B8").Value = "Update variable"
$params.Range("C8").Value = "Attribute Name"
$params.Range("D8").Value = "LastVariableUpdateDate"

$params.Range("A9").Value = "Update_Nightly"
$params.Range("B9").Value = "Update Workspace Plane Status"
$params.Range("C9").Value = "Entity Type Name"
$params.Range("D9").Value = "Workspace_Plan"

$params.Range("A10").Value = "Update_Nightly"
$params.Range("B10").Value = "Update Workspace Plane Status"
$params.Range("C10").Value = "Attribute Name"
$params.Range("D10").Value = "Status"

$params.Range("A11").Value = "Update_Nightly"
$params.Range("B11").Value = "Update Workspace Plane Table"
$params.Range("C11").Value = "formula"
$params.Range("D11").Value = "LIB_EWS.UpdateWorkpalceTableAll();"

$params.Range("A12").Value = "Update_Nightly"
$params.Range("B12").Value = "Update Output table"
$params.Range("C12").Value = "Entity Type Name"
$params.Range("D12").Value = "Analysis_Unit"

$params.Range("A13").Value = "Update_Nightly"
$params.Range("B13").Value = "Update Output table"
$params.Range("C13").Value = "Attribute Name"
$params.Range("D13").Value = "LastUpdateOutput"

$params.Range("A14").Value = "Update_Nightly"
$params.Range("B14").Value = "Update Customer Table Output"
$params.Range("C14").Value = "formula"
$params.Range("D14").Value = "LIB_EWS.UpdateCustomerUnitOutput();"

# Drop the per-cell style overrides on the data rows (A2:E14).
$params.Range("A2:E14").Style = "Normal"

# The "Value" column got noticeably wider to fit the new, longer
# formula text.
$params.Columns.Item(4).ColumnWidth = 42

$params.Activate()
$params.Range("B6").Select()
